$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the requirement rows so that the "LOM3246" entry (Indicacao de
# Conjunto) comes first, followed by the existing Fisica IV and Introducao a
# Ciencia dos Materiais requirements.
$req1 = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"
$req2 = "LOB1021 -  Física IV  (Requisito)`n"
$req3 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"

$ws.Range("B23").Value = $req1
$ws.Range("C23").Value = $req1

$ws.Range("B24").Value = $req2
$ws.Range("C24").Value = $req2

$ws.Range("B25").Value = $req3
$ws.Range("C25").Value = $req3
